$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tactics")

$ws.Range("B4").Value = "(SPD)"
$ws.Range("B7").Value = "(ATK)"
$ws.Range("B8").Value = "(LDR)"
$ws.Range("B3").Value = "(UTL)(DEF)"
$ws.Range("B2").Value = "(INT)(ATK)(SPD)(UTL)"
$ws.Range("B5").Value = "(SPD)(LDR)(ATK)"
$ws.Range("B6").Value = "(SPD)"

$ws.Range("B6").Select()
